$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) column updates ---
# These source values are plain strings in the workbook (t="inlineStr"),
# several of which look like numbers (e.g. "1.00", "20.70", "0.0900").
# A bare .Value assignment lets Excel "smart type" them into numeric
# doubles, which would silently drop the trailing/insignificant zeros.
# Force each target cell to Text format first so the literal string is
# preserved exactly, then restore the default (General/Normal) style so
# no stray number-format is left behind on the cell.
$priceUpdates = [ordered]@{
    'D2' = '37.129.86'
    'D3' = '2.079.54'
    'D4' = '1.00'
    'D5' = '249.89'
    'D6' = '0.649'
    'D8' = '49.84'
    'D9' = '60.47'
    'D11' = '0.0736'
    'D13' = '15.12'
    'D14' = '2.386.20'
    'D15' = '0.826'
    'D16' = '2.086.28'
    'D17' = '5.04'
    'D18' = '37.068.91'
    'D19' = '71.85'
    'D20' = '0.0₃0817'
    'D21' = '13.13'
    'D22' = '237.77'
    'D25' = '2.44'
    'D26' = '168.52'
    'D27' = '9.29'
    'D28' = '20.67'
    'D31' = '1.07'
    'D32' = '4.45'
    'D33' = '0.0602'
    'D34' = '20.70'
    'D35' = '0.0900'
    'D37' = '1.82'
    'D38' = '2.25'
    'D39' = '4.03'
    'D41' = '0.0221'
    'D42' = '17.41'
    'D44' = '97.07'
    'D45' = '2.76'
    'D46' = '0.0875'
    'D47' = '2.97'
    'D48' = '1.301.81'
    'D49' = '6.86'
    'D50' = '2.268.63'
    'D51' = '2.25'
}
foreach ($ref in $priceUpdates.Keys) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $priceUpdates[$ref]
    $ws.Range($ref).Style = "Normal"
}

# --- Coin name (B), link (C) and volume (E) column updates ---
$otherUpdates = [ordered]@{
    'E2' = '  -1.65%  '
    'E3' = '  +6.94%  '
    'E4' = '  +0.43%  '
    'E5' = '  -0.68%  '
    'E6' = '  -6.60%  '
    'E7' = '  +0.09%  '
    'E8' = '  +2.54%  '
    'E9' = '  +2.57%  '
    'E10' = '  -4.52%  '
    'E11' = '  -4.61%  '
    'E12' = '  +4.29%  '
    'E13' = '  -5.12%  '
    'E14' = '  +7.23%  '
    'E15' = '  -2.15%  '
    'E16' = '  +7.41%  '
    'E17' = '  -3.06%  '
    'E18' = '  -1.82%  '
    'E19' = '  -5.22%  '
    'E20' = '  -5.73%  '
    'E21' = '  -4.55%  '
    'E22' = '  -6.87%  '
    'E23' = '  -1.35%  '
    'E24' = '  +0.21%  '
    'E25' = '  -3.41%  '
    'E26' = '  -0.83%  '
    'E27' = '  +3.37%  '
    'E28' = '  +7.98%  '
    'E29' = '  -7.04%  '
    'E30' = '  -6.43%  '
    'B31' = 'ImmutableX'
    'C31' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'E31' = '  +18.12%  '
    'B32' = 'Filecoin'
    'C32' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'E32' = '  -3.83%  '
    'B33' = 'Hedera'
    'C33' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'E33' = '  -2.79%  '
    'B34' = 'Gas'
    'C34' = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
    'E34' = '  +5.18%  '
    'E35' = '  -2.75%  '
    'E36' = '  +0.24%  '
    'B37' = 'WEMIXToken'
    'C37' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'E37' = '  -3.47%  '
    'B38' = 'LidoDAOToken'
    'C38' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'E38' = '  +11.97%  '
    'E39' = '  -7.97%  '
    'E40' = '  -10.11%  '
    'B41' = 'VeChain'
    'C41' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E41' = '  -3.39%  '
    'B42' = 'InjectiveProtocol'
    'C42' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'E42' = '  -1.15%  '
    'E43' = '  +1.56%  '
    'E44' = '  -8.11%  '
    'E45' = '  -4.30%  '
    'E46' = '  +3.04%  '
    'E47' = '  +5.07%  '
    'E48' = '  -4.32%  '
    'E49' = '  +6.18%  '
    'E50' = '  +7.30%  '
    'E51' = '  -8.14%  '
}
foreach ($ref in $otherUpdates.Keys) {
    $ws.Range($ref).Value = $otherUpdates[$ref]
}
